# Fixed bug and added new features
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bug fix: Elmar Qarayev's order for 1440 was marked "Accepted" but should be "Rejected"
$ws.Range("D2").Value = "Rejected"

# New feature: add a new order row for Elmar Qara
$ws.Range("A8").Value = "Elmar Qara"
$ws.Range("B8").Value = "elmarqarayev69@gmail.com"
$ws.Range("C8").Value = 36
$ws.Range("D8").Value = "Pending"
